# Change_of_Schedule_Form: inject {{mustache}} placeholder tags next to
# each fill-in-the-blank label so the template can be mail-merged.
#
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
#   Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)

$d = $word.ActiveDocument

function Replace-Once([string]$find, [string]$replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

Replace-Once "DATE SENT TO SCHEDULER" "DATE SENT TO SCHEDULER {{date_sent}}"
Replace-Once " SPRING " " SPRING {{term_spring}} "
Replace-Once "_______ DIVISION CHAIR" "{{division_chair}}  (DIVISION CHAIR)"
Replace-Once "DATE PROCESSED" "DATE PROCESSED {{date_processed}}"
Replace-Once " SUMMER" " SUMMER {{term_summer}}"
Replace-Once "_______ AREA DEAN" "{{area_dean}}  (AREA DEAN)"
Replace-Once " FALL " " FALL {{term_fall}} "
Replace-Once "VISIBLE IN CLASS SEARCH?" "VISIBLE IN CLASS SEARCH?  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}  YES {{visible_yes}}   NO {{visible_no}}"
Replace-Once "______ LECTURE HOURS" "{{lecture_hours}} LECTURE HOURS"
Replace-Once "______ LAB HOURS" "{{lab_hours}} LAB HOURS"
Replace-Once "______ ACTIVITY HOURS" "{{activity_hours}} ACTIVITY HOURS"
Replace-Once "_______ SEMESTER LECTURE HOURS" "{{sem_lect}} SEMESTER LECTURE HOURS"
Replace-Once "_______ SEMESTER LAB HOURS" "{{sem_lab}} SEMESTER LAB HOURS"
Replace-Once "_______ SEMESTER ACTIVITY HOURS" "{{sem_act}} SEMESTER ACTIVITY HOURS"
Replace-Once "_______ SICK LEAVE HOURS" "{{sick_leave}} SICK LEAVE HOURS"

Write-Output "replacements applied"
